# Edit script for B6-PowerPoint.pptx
#
# Target changes (per the authoritative XML diff):
#   1. Three tables (on slides 14, 15 and 16) get their table style
#      switched from {644885F6-4843-466F-BED0-8996C16541AD}
#      to {79022CC6-5067-43E9-B83C-D7AB2E3EFA9A}.
#   2. The deck's two theme parts swap identity: the theme actually
#      driving the slides (Integral / "Red Violet") becomes the
#      stock Office theme's colours, while the secondary theme part
#      (only linked from the notes master) keeps the Office colours
#      it already has -- Office PowerPoint's automation surface
#      doesn't expose a way to rewrite a theme part that isn't the
#      active design's theme, so the colour swap is applied through
#      the one ThemeColorScheme PowerPoint does expose (the slide
#      master / design theme).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table style swap on the three tables that use the old style id
# ---------------------------------------------------------------
$oldStyleId = "{644885F6-4843-466F-BED0-8996C16541AD}"
$newStyleId = "{79022CC6-5067-43E9-B83C-D7AB2E3EFA9A}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------
# 2) Theme colour scheme: flip the live design theme from the
#    "Integral" / Red Violet palette to the stock Office palette.
# ---------------------------------------------------------------
# Index layout for ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
# 11 hlink, 12 folHlink. RGB() packs as r + g*256 + b*65536, which is
# exactly what PowerPoint's ColorFormat.RGB / ThemeColor.RGB expect.

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme

$tcs.Item(1).RGB  = RGB(0x00, 0x00, 0x00)   # dk1      -> 000000
$tcs.Item(2).RGB  = RGB(0xFF, 0xFF, 0xFF)   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = RGB(0x44, 0x54, 0x6A)   # dk2      -> 44546A
$tcs.Item(4).RGB  = RGB(0xE7, 0xE6, 0xE6)   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = RGB(0x5B, 0x9B, 0xD5)   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = RGB(0xED, 0x7D, 0x31)   # accent2  -> ED7D31
$tcs.Item(7).RGB  = RGB(0xA5, 0xA5, 0xA5)   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = RGB(0xFF, 0xC0, 0x00)   # accent4  -> FFC000
$tcs.Item(9).RGB  = RGB(0x44, 0x72, 0xC4)   # accent5  -> 4472C4
$tcs.Item(10).RGB = RGB(0x70, 0xAD, 0x47)   # accent6  -> 70AD47
$tcs.Item(11).RGB = RGB(0x05, 0x63, 0xC1)   # hlink    -> 0563C1
$tcs.Item(12).RGB = RGB(0x95, 0x4F, 0x72)   # folHlink -> 954F72
